# Change "MP3 module" (single row) to a DFPlayer Mini MP3 Player wired with
# separate TX/RX pins ("MP3 TX" / "MP3 RX"), shifting the stepper-driver
# labels down to make room.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pin")

# New TX/RX pins for the MP3 player take the rows previously used for the
# stepper driver's "enable"/"pulse" labels (rows 6 & 7); the stepper labels
# move down into the rows vacated by the old single "MP3 module" row.
# (leading apostrophe forces text / matches the quote-prefixed style used
# by the other "-" placeholder cells in this column)
$ws.Range("H6").Value = "'MP3 TX"
$ws.Range("H7").Value = "'MP3 RX"

$ws.Range("H12").Value = "Stepper pulse"
$ws.Range("H13").Value = "Stepper direction"

# The old "MP3 module" row is now blank.
$ws.Range("H9").Value = ""

$ws.Range("H9").Select()
